$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep a Text format so numeric-looking strings
# (e.g. "119.60", "1.00") are not coerced into numbers, matching the
# inline-string representation used in the source workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.685.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.273.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "119.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "265.97"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.22"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.55%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.49"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.908"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.614.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.277.57"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.659.08"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.40"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.52"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.99"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.98"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.63"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0919"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.78"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0388"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.23"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +13.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.58"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.55"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.78%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.82"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.79"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.80"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.30"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +41.94%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.55"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.96"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.69%  "
